# -----------------------------------------------------------------------
# Commit: "Removed sys appends for dev"
#
# 1) Sheets rId3 / rId5 swap their display names:
#      rId3: "Criterion 3, Air Speed 0.1" -> "Criterion 1, Air Speed 0.1"
#      rId5: "Criterion 1, Air Speed 0.1" -> "Criterion 3, Air Speed 0.1"
# 2) The data that used to live on rId3 (the "Criterion 3" sheet, which only
#    ever held placeholder zeros) now lives on rId5, and vice-versa: the
#    real per-room Absolute/Relative change numbers move from the old
#    "Criterion 1" sheet onto what is now the "Criterion 1" sheet (rId3),
#    and the placeholder zeros move onto what is now the "Criterion 3"
#    sheet (rId5).
# 3) The "readme" manifest table's column order changes from
#    index/sheet_name/Author/JobNo/Date to index/Author/sheet_name/Date/JobNo,
#    and its data rows are rewritten to match.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Step 1: swap the two worksheet names (rId3 <-> rId5) ---------------
$ws3 = $wb.Worksheets.Item("Criterion 3, Air Speed 0.1")
$ws5 = $wb.Worksheets.Item("Criterion 1, Air Speed 0.1")
$ws3.Name = "__tmp_swap_sheet_name__"
$ws5.Name = "Criterion 3, Air Speed 0.1"
$ws3.Name = "Criterion 1, Air Speed 0.1"

# --- Step 2: re-point the two worksheet variables at the renamed sheets -
$ws3 = $wb.Worksheets.Item("Criterion 1, Air Speed 0.1")
$ws5 = $wb.Worksheets.Item("Criterion 3, Air Speed 0.1")

# --- Step 3: update the column header labels (E1/F1) on each sheet ------
$ws3.Range("E1").Value = "Criterion 1 Absolute Change"
$ws3.Range("F1").Value = "Criterion 1 Relative Change (%)"
$ws5.Range("E1").Value = "Criterion 3 Absolute Change"
$ws5.Range("F1").Value = "Criterion 3 Relative Change (%)"

# --- Step 4: rewrite the data payload (C2:F28) that now belongs on the
# "Criterion 1" sheet (was on the "Criterion 3" sheet before the swap) ---

$ws3.Range("C2").Value = 2.5
$ws3.Range("D2").Value = 2.48
$ws3.Range("E2").Value = -0.02000000000000002
$ws3.Range("F2").Value = -0.8000000000000007
$ws3.Range("C3").Value = 7.8
$ws3.Range("D3").Value = 7.89
$ws3.Range("E3").Value = 0.08999999999999986
$ws3.Range("F3").Value = 1.153846153846152
$ws3.Range("C4").Value = 8.1
$ws3.Range("D4").Value = 8.039999999999999
$ws3.Range("E4").Value = -0.0600000000000005
$ws3.Range("F4").Value = -0.7407407407407469
$ws3.Range("C5").Value = 0
$ws3.Range("D5").Value = 0
$ws3.Range("E5").Value = 0
$ws3.Range("C6").Value = 0
$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 0
$ws3.Range("C7").Value = 0
$ws3.Range("D7").Value = 0
$ws3.Range("E7").Value = 0
$ws3.Range("C8").Value = 0
$ws3.Range("D8").Value = 0
$ws3.Range("E8").Value = 0
$ws3.Range("C9").Value = 10.6
$ws3.Range("D9").Value = 10.59
$ws3.Range("E9").Value = -0.009999999999999787
$ws3.Range("F9").Value = -0.09433962264150743
$ws3.Range("C10").Value = 26
$ws3.Range("D10").Value = 26.04
$ws3.Range("E10").Value = 0.03999999999999915
$ws3.Range("F10").Value = 0.1538461538461506
$ws3.Range("C11").Value = 11.6
$ws3.Range("D11").Value = 11.57
$ws3.Range("E11").Value = -0.02999999999999936
$ws3.Range("F11").Value = -0.2586206896551669
$ws3.Range("C12").Value = 26.3
$ws3.Range("D12").Value = 26.4
$ws3.Range("E12").Value = 0.09999999999999787
$ws3.Range("F12").Value = 0.3802281368821211
$ws3.Range("C13").Value = 0
$ws3.Range("D13").Value = 0
$ws3.Range("E13").Value = 0
$ws3.Range("C14").Value = 0.6
$ws3.Range("D14").Value = 0.6
$ws3.Range("E14").Value = 0
$ws3.Range("F14").Value = 0
$ws3.Range("C15").Value = 0.1
$ws3.Range("D15").Value = 0.14
$ws3.Range("E15").Value = 0.04000000000000001
$ws3.Range("F15").Value = 40.00000000000001
$ws3.Range("C16").Value = 0
$ws3.Range("D16").Value = 0
$ws3.Range("E16").Value = 0
$ws3.Range("C17").Value = 12.9
$ws3.Range("D17").Value = 12.94
$ws3.Range("E17").Value = 0.03999999999999915
$ws3.Range("F17").Value = 0.3100775193798383
$ws3.Range("C18").Value = 28.3
$ws3.Range("D18").Value = 28.26
$ws3.Range("E18").Value = -0.03999999999999915
$ws3.Range("F18").Value = -0.1413427561837426
$ws3.Range("C19").Value = 13.4
$ws3.Range("D19").Value = 13.45
$ws3.Range("E19").Value = 0.04999999999999893
$ws3.Range("F19").Value = 0.373134328358201
$ws3.Range("C20").Value = 28.6
$ws3.Range("D20").Value = 28.46
$ws3.Range("E20").Value = -0.1400000000000006
$ws3.Range("F20").Value = -0.4895104895104915
$ws3.Range("C21").Value = 0
$ws3.Range("D21").Value = 0.03
$ws3.Range("E21").Value = 0.03
$ws3.Range("F21").Value = "inf"
$ws3.Range("C22").Value = 1.6
$ws3.Range("D22").Value = 1.58
$ws3.Range("E22").Value = -0.02000000000000002
$ws3.Range("F22").Value = -1.250000000000001
$ws3.Range("C23").Value = 0.4
$ws3.Range("D23").Value = 0.35
$ws3.Range("E23").Value = -0.05000000000000004
$ws3.Range("F23").Value = -12.50000000000001
$ws3.Range("C24").Value = 0
$ws3.Range("D24").Value = 0.03
$ws3.Range("E24").Value = 0.03
$ws3.Range("F24").Value = "inf"
$ws3.Range("C25").Value = 5.6
$ws3.Range("D25").Value = 5.61
$ws3.Range("E25").Value = 0.01000000000000068
$ws3.Range("F25").Value = 0.1785714285714406
$ws3.Range("C26").Value = 17.7
$ws3.Range("D26").Value = 17.8
$ws3.Range("E26").Value = 0.1000000000000014
$ws3.Range("F26").Value = 0.5649717514124375
$ws3.Range("C27").Value = 5.2
$ws3.Range("D27").Value = 5.23
$ws3.Range("E27").Value = 0.03000000000000025
$ws3.Range("F27").Value = 0.5769230769230818
$ws3.Range("C28").Value = 17.7
$ws3.Range("D28").Value = 17.75
$ws3.Range("E28").Value = 0.05000000000000071
$ws3.Range("F28").Value = 0.2824858757062187


# --- Step 5: rewrite the data payload (C2:F28) that now belongs on the
# "Criterion 3" sheet (was on the "Criterion 1" sheet before the swap) ---

$ws5.Range("C2").Value = 2
$ws5.Range("D2").Value = 2
$ws5.Range("E2").Value = 0
$ws5.Range("F2").Value = 0
$ws5.Range("C3").Value = 2
$ws5.Range("D3").Value = 2
$ws5.Range("E3").Value = 0
$ws5.Range("F3").Value = 0
$ws5.Range("C4").Value = 2
$ws5.Range("D4").Value = 2
$ws5.Range("E4").Value = 0
$ws5.Range("F4").Value = 0
$ws5.Range("C5").Value = 0
$ws5.Range("D5").Value = 0
$ws5.Range("E5").Value = 0
$ws5.Range("C6").Value = 0
$ws5.Range("D6").Value = 0
$ws5.Range("E6").Value = 0
$ws5.Range("C7").Value = 0
$ws5.Range("D7").Value = 0
$ws5.Range("E7").Value = 0
$ws5.Range("C8").Value = 0
$ws5.Range("D8").Value = 0
$ws5.Range("E8").Value = 0
$ws5.Range("C9").Value = 3
$ws5.Range("D9").Value = 3
$ws5.Range("E9").Value = 0
$ws5.Range("F9").Value = 0
$ws5.Range("C10").Value = 3
$ws5.Range("D10").Value = 3
$ws5.Range("E10").Value = 0
$ws5.Range("F10").Value = 0
$ws5.Range("C11").Value = 3
$ws5.Range("D11").Value = 3
$ws5.Range("E11").Value = 0
$ws5.Range("F11").Value = 0
$ws5.Range("C12").Value = 3
$ws5.Range("D12").Value = 3
$ws5.Range("E12").Value = 0
$ws5.Range("F12").Value = 0
$ws5.Range("C13").Value = 0
$ws5.Range("D13").Value = 0
$ws5.Range("E13").Value = 0
$ws5.Range("C14").Value = 1
$ws5.Range("D14").Value = 1
$ws5.Range("E14").Value = 0
$ws5.Range("F14").Value = 0
$ws5.Range("C15").Value = 1
$ws5.Range("D15").Value = 1
$ws5.Range("E15").Value = 0
$ws5.Range("F15").Value = 0
$ws5.Range("C16").Value = 0
$ws5.Range("D16").Value = 0
$ws5.Range("E16").Value = 0
$ws5.Range("C17").Value = 3
$ws5.Range("D17").Value = 3
$ws5.Range("E17").Value = 0
$ws5.Range("F17").Value = 0
$ws5.Range("C18").Value = 3
$ws5.Range("D18").Value = 3
$ws5.Range("E18").Value = 0
$ws5.Range("F18").Value = 0
$ws5.Range("C19").Value = 3
$ws5.Range("D19").Value = 3
$ws5.Range("E19").Value = 0
$ws5.Range("F19").Value = 0
$ws5.Range("C20").Value = 3
$ws5.Range("D20").Value = 3
$ws5.Range("E20").Value = 0
$ws5.Range("F20").Value = 0
$ws5.Range("C21").Value = 1
$ws5.Range("D21").Value = 1
$ws5.Range("E21").Value = 0
$ws5.Range("F21").Value = 0
$ws5.Range("C22").Value = 1
$ws5.Range("D22").Value = 1
$ws5.Range("E22").Value = 0
$ws5.Range("F22").Value = 0
$ws5.Range("C23").Value = 1
$ws5.Range("D23").Value = 1
$ws5.Range("E23").Value = 0
$ws5.Range("F23").Value = 0
$ws5.Range("C24").Value = 1
$ws5.Range("D24").Value = 1
$ws5.Range("E24").Value = 0
$ws5.Range("F24").Value = 0
$ws5.Range("C25").Value = 2
$ws5.Range("D25").Value = 2
$ws5.Range("E25").Value = 0
$ws5.Range("F25").Value = 0
$ws5.Range("C26").Value = 2
$ws5.Range("D26").Value = 2
$ws5.Range("E26").Value = 0
$ws5.Range("F26").Value = 0
$ws5.Range("C27").Value = 2
$ws5.Range("D27").Value = 2
$ws5.Range("E27").Value = 0
$ws5.Range("F27").Value = 0
$ws5.Range("C28").Value = 2
$ws5.Range("D28").Value = 2
$ws5.Range("E28").Value = 0
$ws5.Range("F28").Value = 0


# --- Step 6: reorder the "readme" manifest columns and rewrite its data -
$ws1 = $wb.Worksheets.Item("readme")

$ws1.Range("B1").Value = "Author"
$ws1.Range("C1").Value = "sheet_name"
$ws1.Range("D1").Value = "Date"
$ws1.Range("E1").Value = "JobNo"
$ws1.Range("B2").Value = "jovyan"
$ws1.Range("C2").Value = "Criteria Failing, Air Speed 0.1"
$ws1.Range("D2").Value = "20220308"
$ws1.Range("E2").Value = "/c/e"
$ws1.Range("B3").Value = "jovyan"
$ws1.Range("C3").Value = "Criterion 1, Air Speed 0.1"
$ws1.Range("D3").Value = "20220308"
$ws1.Range("E3").Value = "/c/e"
$ws1.Range("B4").Value = "jovyan"
$ws1.Range("C4").Value = "Criterion 2, Air Speed 0.1"
$ws1.Range("D4").Value = "20220308"
$ws1.Range("E4").Value = "/c/e"
$ws1.Range("B5").Value = "jovyan"
$ws1.Range("C5").Value = "Criterion 3, Air Speed 0.1"
$ws1.Range("D5").Value = "20220308"
$ws1.Range("E5").Value = "/c/e"


Write-Output "edit applied"
